$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained one new data row: a new record was inserted at row 25
# ("Castle Brite" / "Provincia de San Felipe de Aconcagua", 2021-12-09),
# pushing the existing rows 25-42 down to rows 26-43.

# 1) Shift existing rows 25..42 down to 26..43 (iterate bottom-up so we
#    never overwrite a source row before it has been read).
for ($r = 42; $r -ge 25; $r--) {
    $target = $r + 1

    # Make sure the destination date cell (column D) carries the same
    # date number-format as the source, since it's a brand-new cell when
    # target = 43.
    $ws.Cells.Item($target, 4).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat

    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($target, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Write the new record into row 25.
$ws.Cells.Item(25, 1).Value  = 10
$ws.Cells.Item(25, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value  = "La Araucanía"
$ws.Cells.Item(25, 4).Value  = 44539
$ws.Cells.Item(25, 5).Value  = 9
$ws.Cells.Item(25, 6).Value  = "Fruta"
$ws.Cells.Item(25, 7).Value  = 100103
$ws.Cells.Item(25, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(25, 9).Value  = 100103003
$ws.Cells.Item(25, 10).Value = "Damasco"
$ws.Cells.Item(25, 11).Value = "Castle Brite"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 75
$ws.Cells.Item(25, 14).Value = 18000
$ws.Cells.Item(25, 15).Value = 20000
$ws.Cells.Item(25, 16).Value = 19067
$ws.Cells.Item(25, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(25, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(25, 19).Value = 1059
$ws.Cells.Item(25, 20).Value = 18
